$d = $word.ActiveDocument
$q = [char]0x2019

function FR($old, $new) {
    $f = $d.Content.Find
    $f.ClearFormatting()
    $f.Replacement.ClearFormatting()
    $ok = $f.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "FIND/REPLACE FAILED:" $old
    }
}

# 1. "McDonald's international menus, which could be featured in the" -> "... menus that should be featured in the"
FR "McDonald${q}s international menus, which could be featured in the " "McDonald${q}s international menus that should be featured in the "

# 2. "sought a new sales approach ... ethnicities.  To this end, McDonald's decided to initiate" -> "sought a new approach that would appeal to these Americans and thus, decided to initiate"
FR "McDonald${q}s sought a new sales approach by marketing to the country${q}s diverse cultures and ethnicities.  To this end, McDonald${q}s decided to initiate" "McDonald${q}s sought a new approach that would appeal to these Americans and thus, decided to initiate"

# 3. Cascade Range sentence simplified to "Its climate allows for"
FR "coast of the United States.  The Cascade Range vertically splits the state, forming two contrasting climates: East of the Cascades, and West of the Cascades.  This diversity allows for " "coast of the United States.  Its climate allows for "

# 4. "To the West of the Cascades, the coastal region yields" -> "The coastal region of Oregon also yields"
FR "To the West of the Cascades, the coastal region yields" "The coastal region of Oregon also yields"

# 5. "asked about their eating preferences.  Results from the" -> "asked about their eating preferences and opinion on McDonald's.  Results from the"
FR " asked about their eating preferences.  Results from the " " asked about their eating preferences and opinion on McDonald${q}s.  Results from the "

# 6. "Italy's Spinach and Parmesan Cheese Nuggets" -> "Italy's Spinach & Parmesan Cheese Nuggets"
FR "Italy${q}s Spinach and Parmesan Cheese Nuggets" "Italy${q}s Spinach & Parmesan Cheese Nuggets"

# 6b. Word's "_GoBack" bookmark tracks the most recent edit location; move it to just
# after the freshly-typed "&" (mirrors what Word does after an in-place edit there).
$fGoBack = $d.Content.Find
$fGoBack.ClearFormatting()
$okGoBack = $fGoBack.Execute("Italy${q}s Spinach &", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($okGoBack) {
    $goBackPos = $fGoBack.Parent.End
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos))
}

# 7. "can easily be provided, fresh, by" -> "can easily and freshly provided by"
FR "can easily be provided, fresh, by" "can easily and freshly provided by"

Write-Host "All replacements attempted."
